$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.753.24"
$ws.Range("E2").Value = "  +3.52%  "

$ws.Range("D3").Value = "1.923.37"
$ws.Range("E3").Value = "  +2.31%  "

$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -2.15%  "

$ws.Range("D5").Value = "335.01"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -1.53%  "

$ws.Range("D7").Value = "0.4673"
$ws.Range("E7").Value = "  +1.69%  "

$ws.Range("D8").Value = "0.4145"
$ws.Range("E8").Value = "  +4.01%  "

$ws.Range("D9").Value = "48.32"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").Value = "0.08053"
$ws.Range("E10").Value = "  +0.99%  "

$ws.Range("D11").Value = "1.019"
$ws.Range("E11").Value = "  +2.53%  "

$ws.Range("D12").Value = "22.32"
$ws.Range("E12").Value = "  +2.51%  "

$ws.Range("D13").Value = "1.909.17"
$ws.Range("E13").Value = "  +1.23%  "

$ws.Range("D14").Value = "6.023"
$ws.Range("E14").Value = "  +1.31%  "

$ws.Range("D15").Value = "7.207"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").Value = "89.93"
$ws.Range("E16").Value = "  +0.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -2.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001037"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").Value = "0.06592"
$ws.Range("E19").Value = "  -0.60%  "

$ws.Range("D20").Value = "17.86"
$ws.Range("E20").Value = "  +3.29%  "

$ws.Range("D21").Value = "0.9991"
$ws.Range("E21").Value = "  -1.20%  "

$ws.Range("D22").Value = "29.704.27"
$ws.Range("E22").Value = "  +2.93%  "

$ws.Range("D23").Value = "5.562"
$ws.Range("E23").Value = "  +2.31%  "

$ws.Range("D24").Value = "11.67"
$ws.Range("E24").Value = "  +7.43%  "

$ws.Range("D25").Value = "2.198"
$ws.Range("E25").Value = "  -3.01%  "

$ws.Range("D26").Value = "2.139.51"
$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("D27").Value = "157.31"
$ws.Range("E27").Value = "  -1.27%  "

$ws.Range("D28").Value = "19.95"
$ws.Range("E28").Value = "  +1.71%  "

$ws.Range("D29").Value = "2.151"
$ws.Range("E29").Value = "  +2.98%  "

$ws.Range("D30").Value = "5.743"
$ws.Range("E30").Value = "  +6.71%  "

$ws.Range("D31").Value = "117.91"
$ws.Range("E31").Value = "  -1.33%  "

$ws.Range("E32").Value = "  +8.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09470"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").Value = "1.437"
$ws.Range("E34").Value = "  +1.49%  "

$ws.Range("D35").Value = "5.438"
$ws.Range("E35").Value = "  +2.30%  "

$ws.Range("D36").Value = "3.523"
$ws.Range("E36").Value = "  -2.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06160"
$ws.Range("E37").Value = "  +1.30%  "

$ws.Range("E38").Value = "  +1.17%  "

$ws.Range("D39").Value = "8.478"
$ws.Range("E39").Value = "  +3.05%  "

$ws.Range("D40").Value = "1.184"
$ws.Range("E40").Value = "  +1.95%  "

$ws.Range("D41").Value = "0.5922"
$ws.Range("E41").Value = "  +1.38%  "

$ws.Range("D42").Value = "0.1849"
$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("D43").Value = "10.27"
$ws.Range("E43").Value = "  +0.52%  "

$ws.Range("D44").Value = "2.352"
$ws.Range("E44").Value = "  +2.79%  "

$ws.Range("D45").Value = "1.239"
$ws.Range("E45").Value = "  -1.24%  "

$ws.Range("D46").Value = "0.07525"
$ws.Range("E46").Value = "  +3.41%  "

$ws.Range("D47").Value = "0.5604"
$ws.Range("E47").Value = "  +1.97%  "

$ws.Range("D48").Value = "12.29"
$ws.Range("E48").Value = "  +0.27%  "

$ws.Range("E49").Value = "  +1.79%  "

$ws.Range("D50").Value = "113.02"
$ws.Range("E50").Value = "  +1.87%  "

$ws.Range("D51").Value = "0.3001"
$ws.Range("E51").Value = "  +11.69%  "
